# Time Planning.xlsx update: move three "newly implemented" cards from the
# "Text Cards to implement" column (C) to the "Finished Cards" column (D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Militia (C18) -> D37, Spy (C23) -> D38, Money Lender (C21) -> D39
$ws.Range("D37").Value = $ws.Range("C18").Value2
$ws.Range("D38").Value = $ws.Range("C23").Value2
$ws.Range("D39").Value = $ws.Range("C21").Value2

# Clear the now-vacated source cells in column C
$ws.Range("C18").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("C23").ClearContents()

# Update the view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("D39").Select()

# Update workbook window position
$excel.ActiveWindow.Left = 7188
